$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B/C (text: coin name / link) are plain text and need no
# special handling. Cells in columns D/E (price / volume%) are numeric-looking
# strings in the source data; Excel auto-converts those to real numbers on
# assignment (dropping trailing zeros / percent formatting), so we force the
# cell format to Text ("@") before assigning, preserving exact source text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.20%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.045"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.19%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07599"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.41%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.598"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.14%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.448"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.33%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9056"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1023"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.54%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1756"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.43%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09046"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.13%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04228"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.51%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1053"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.32%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001259"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005830"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.95%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.355"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.46%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.270"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.12%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3268"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.781"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.66%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1357"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.11%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.53%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.49%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004064"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.21%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.97%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003012"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.41%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02395"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.51%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05165"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007798"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.85%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1304"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.84%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007071"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.90%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001924"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.74%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008266"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.18%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3350"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.15%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006363"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.57%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004408"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.75%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.006616"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "96.85%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
